# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.252.57"
$ws.Range("E2").Value = "  +2.81%  "

$ws.Range("D3").Value = "1.718.59"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.75"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("E7").Value = "  -1.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2621"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06193"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").Value = "1.717.67"
$ws.Range("E10").Value = "  +3.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07073"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.33"
$ws.Range("E12").Value = "  +3.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5968"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.420"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.16"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "26.262.08"
$ws.Range("E18").Value = "  +2.89%  "

$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D21").Value = "1.937.47"
$ws.Range("E21").Value = "  +3.38%  "

$ws.Range("E22").Value = "  +2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.724"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.281"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.85"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.403"
$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.761"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.28"
$ws.Range("E29").Value = "  +1.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.965"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07754"
$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04454"
$ws.Range("E33").Value = "  +4.87%  "

$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9999"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.616"
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9739"
$ws.Range("E36").Value = "  +2.51%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6183"
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9242"
$ws.Range("E38").Value = "  +7.21%  "

$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "113.71"
$ws.Range("E39").Value = "  +17.27%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.412"
$ws.Range("E40").Value = "  -7.28%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.922"
$ws.Range("E41").Value = "  +3.53%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01484"
$ws.Range("E43").Value = "  +1.17%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.605"
$ws.Range("E44").Value = "  +15.65%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3822"
$ws.Range("E45").Value = "  +1.47%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1181"
$ws.Range("E46").Value = "  +5.17%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.288"
$ws.Range("E47").Value = "  +1.33%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05273"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.818"
$ws.Range("E49").Value = "  +6.29%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.30"
$ws.Range("E50").Value = "  +1.86%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3372"
$ws.Range("E51").Value = "  +1.10%  "
